# Auto-committed on 2023/07/21 週五 17:18:54.18
#
# Applies the tracked edits to the "DBD" worksheet of CdCode.xlsx:
#   - a couple of quantity values were bumped (E11, E15, E16)
#   - the formatting of row 17 (and the leading "#" column for rows
#     18-21) was switched from a set of near-duplicate cell styles
#     (which had an extra, unused "applyFill" flag) back onto the
#     already-existing, equivalent styles used by the rest of the
#     table, by copying formats across with PasteSpecial
#   - the active sheet's selection moved to H20

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# --- value edits -----------------------------------------------------
$ws.Range("E11").Value = 30
$ws.Range("E15").Value = 2
$ws.Range("E16").Value = 2

# --- style clean-up ----------------------------------------------------
# Row 17's A/D/F/G cells (and the "#" cells of rows 18-21) used a style
# that only differed from the style already used elsewhere in the "#"
# column (e.g. A9) by a stray, no-op applyFill flag -- reuse that style.
$ws.Range("A9").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("F17").PasteSpecial(-4122)
$ws.Range("G17").PasteSpecial(-4122)
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A21").PasteSpecial(-4122)

# Row 17's B/C/E cells used a similarly redundant style; G24 already
# carries the de-duplicated equivalent.
$ws.Range("G24").Copy()
$ws.Range("B17:C17").PasteSpecial(-4122)
$ws.Range("E17").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- selection ---------------------------------------------------------
$ws.Range("H20").Select() | Out-Null
